$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2565f4f8b28b1fb60927bb690b4be35636bb3f04/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2565f4f8b28b1fb60927bb690b4be35636bb3f04/e2e/b.md"

$hyperlinkUnderline = 2          # xlUnderlineStyleSingle
$hyperlinkColor = 15570276       # BGR encoding of FF6495ED (matches existing workbook HyperLink style)

function Set-HyperlinkStyle($rng) {
    $rng.Font.Underline = $hyperlinkUnderline
    $rng.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# 1. Update status text on the Overview sheet (report generated for handback)
# ---------------------------------------------------------------------------
$ws1.Range('E2').Value = $newStatus
$ws1.Range('F2').Value = $newStatus
$ws1.Range('E3').Value = $newStatus
$ws1.Range('F3').Value = $newStatus

$ws1.Columns("E:F").AutoFit()
$ws1.Range('E1').ColumnWidth = 29.1
$ws1.Range('F1').ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: record handback of a.md (rows 2 and 3 both describe a.md)
# ---------------------------------------------------------------------------
$ws2.Range('C2').Value = $newStatus
$ws2.Range('C3').Value = $newStatus

$ws2.Hyperlinks.Add($ws2.Range('I2'), $urlA, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
Set-HyperlinkStyle $ws2.Range('I2')
$ws2.Range('J2').Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws2.Range('K2').Value = "2016-08-15 12:31:18"

$ws2.Hyperlinks.Add($ws2.Range('I3'), $urlA, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
Set-HyperlinkStyle $ws2.Range('I3')
$ws2.Range('J3').Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$ws2.Range('K3').Value = "2016-08-15 12:31:18"

$ws2.Columns("C:C").AutoFit()
$ws2.Range('C1').ColumnWidth = 29.1
$ws2.Range('J1').ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# 3. de-de sheet: record handback of a.md (rows 2 and 3 both describe a.md)
# ---------------------------------------------------------------------------
$ws3.Range('C2').Value = $newStatus
$ws3.Range('C3').Value = $newStatus

$ws3.Hyperlinks.Add($ws3.Range('I2'), $urlA, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
Set-HyperlinkStyle $ws3.Range('I2')
$ws3.Range('J2').Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws3.Range('K2').Value = "2016-08-15 12:31:24"

$ws3.Hyperlinks.Add($ws3.Range('I3'), $urlA, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "a.md")
Set-HyperlinkStyle $ws3.Range('I3')
$ws3.Range('J3').Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$ws3.Range('K3').Value = "2016-08-15 12:31:24"

$ws3.Columns("C:C").AutoFit()
$ws3.Range('C1').ColumnWidth = 29.1
$ws3.Range('J1').ColumnWidth = 39.17

Write-Host "Report generated for handback."
